$d = $word.ActiveDocument
$d.Content.Find.Execute("Upon reviewing the application made and upon considering the information provided by the parties, the court requests more information from the applicant.", $true, $false, $false, $false, $false, $true, 1, $false, "<<judgeRecital>>", 2)
